$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Rebuild the whole body story with the final wording / formatting, using
# InsertXML so every paragraph/run comes out exactly as specified (no leftover
# Times New Roman rFonts, no inherited rsid/paraId cruft). The trailing
# <w:sectPr> is left untouched by Content/InsertXML, so it is not included
# here - margins are adjusted separately below via PageSetup.
# ---------------------------------------------------------------------------
$bodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>HỢP ĐỒNG TÍN DỤNG</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Kính gửi khách hàng: </w:t></w:r><w:r><w:t>Đặng Ngọc Mai</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">CCCD: </w:t></w:r><w:r><w:t>33301007320</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Số tiền cấp: </w:t></w:r><w:r><w:t>6,160,000 VND</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Ngày tạo hợp đồng: </w:t></w:r><w:r><w:t>08/10/2025</w:t></w:r></w:p><w:p><w:r><w:t>Sign here:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>---------------------</w:t></w:r></w:p><w:p/>
'@

$d.Content.InsertXML($bodyXml)

# ---------------------------------------------------------------------------
# Page margins: 1418/1134/1418/1985 twips -> 1440/1440/1440/1440 twips
# (72pt / 72pt / 72pt / 72pt)
# ---------------------------------------------------------------------------
$ps = $d.PageSetup
$ps.TopMargin = 72
$ps.RightMargin = 72
$ps.BottomMargin = 72
$ps.LeftMargin = 72

# ---------------------------------------------------------------------------
# Drop the now-unused "TOC 1" paragraph style.
# ---------------------------------------------------------------------------
$d.Styles("TOC 1").Delete()
